$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 19
$ws.Range("H19").Value = 383.4
$ws.Range("J19").Value = 279.64706
$ws.Range("L19").Value = 279.64706
$ws.Range("N19").Value = -629.64706
# row 33
$ws.Range("H33").Value = 499.35715
$ws.Range("I33").Value = 537.4167
$ws.Range("J33").Value = 271
$ws.Range("K33").Value = 537.4167
$ws.Range("L33").Value = 271
$ws.Range("M33").Value = -308.4167
$ws.Range("N33").Value = -729
# row 125
$ws.Range("H125").Value = 2639.04
$ws.Range("I125").Value = 2642
$ws.Range("J125").Value = 2636.7144
$ws.Range("K125").Value = 23778
$ws.Range("L125").Value = 23730.4296
$ws.Range("M125").Value = -21318
$ws.Range("N125").Value = -28650.4296
# row 134
$ws.Range("H134").Value = 38513
$ws.Range("J134").Value = 38513
$ws.Range("L134").Value = 38513
$ws.Range("N134").Value = -48653
# row 137
$ws.Range("H137").Value = 2275.5264
$ws.Range("I137").Value = 3008.3333
$ws.Range("J137").Value = 1937.3077
$ws.Range("K137").Value = 9024.999899999999
$ws.Range("L137").Value = 5811.9231
$ws.Range("M137").Value = -6474.999899999999
$ws.Range("N137").Value = -10911.9231
# row 138
$ws.Range("H138").Value = 2738.7856
$ws.Range("I138").Value = 2687.818
$ws.Range("J138").Value = 2925.6667
$ws.Range("K138").Value = 8063.454000000001
$ws.Range("L138").Value = 8777.000100000001
$ws.Range("M138").Value = -2923.454000000001
$ws.Range("N138").Value = -19057.0001

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 6209.69
$ws.Range("I32").Value = 4887.62
$ws.Range("J32").Value = 11183.19
$ws.Range("K32").Value = 4887.62
$ws.Range("L32").Value = 11183.19
$ws.Range("M32").Value = -4600.62
$ws.Range("N32").Value = -11757.19
# row 96
$ws.Range("H96").Value = 19996.334
$ws.Range("J96").Value = 19996.334
$ws.Range("L96").Value = 19996.334
$ws.Range("N96").Value = -25488.334
# row 102
$ws.Range("H102").Value = 55561492
$ws.Range("I102").Value = 55561492
$ws.Range("K102").Value = 55561492
$ws.Range("M102").Value = -55559870
# row 132
$ws.Range("H132").Value = 2204.9167
$ws.Range("I132").Value = 2030.2759
$ws.Range("J132").Value = 2928.4285
$ws.Range("K132").Value = 6090.8277
$ws.Range("L132").Value = 8785.2855
$ws.Range("M132").Value = -3560.8277
$ws.Range("N132").Value = -13845.2855
# row 141
$ws.Range("H141").Value = 52978
$ws.Range("J141").Value = 52978
$ws.Range("L141").Value = 52978
$ws.Range("N141").Value = -63338

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 75
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
# row 78
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
# row 96
$ws.Range("H96").Value = 20000.5
$ws.Range("J96").Value = 20000.5
$ws.Range("L96").Value = 20000.5
$ws.Range("N96").Value = -25492.5
# row 107
$ws.Range("H107").Value = 1259.2307
$ws.Range("I107").Value = 653.75
$ws.Range("J107").Value = 2228
$ws.Range("K107").Value = 653.75
$ws.Range("L107").Value = 2228
$ws.Range("M107").Value = 1266.25
$ws.Range("N107").Value = -6068
# row 132
$ws.Range("H132").Value = 1794.4791
$ws.Range("I132").Value = 1448.0264
$ws.Range("J132").Value = 3111
$ws.Range("K132").Value = 4344.0792
$ws.Range("L132").Value = 9333
$ws.Range("M132").Value = -1814.0792
$ws.Range("N132").Value = -14393
# row 134
$ws.Range("H134").Value = 11629402
$ws.Range("I134").Value = 1476.7241
$ws.Range("J134").Value = 35715820
$ws.Range("K134").Value = 4430.1723
$ws.Range("L134").Value = 107147460
$ws.Range("M134").Value = -1895.1723
$ws.Range("N134").Value = -107152530
# row 140
$ws.Range("H140").Value = 64000
$ws.Range("J140").Value = 64000
$ws.Range("L140").Value = 64000
$ws.Range("N140").Value = -74360

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 14
$ws.Range("H14").Value = 256.1111
$ws.Range("I14").Value = 256.1111
$ws.Range("K14").Value = 768.3333
$ws.Range("M14").Value = -595.3333
# row 32
$ws.Range("H32").Value = 2436.3635
$ws.Range("J32").Value = 2266.6667
$ws.Range("L32").Value = 6800.000100000001
$ws.Range("N32").Value = -7366.000100000001
# row 122
$ws.Range("H122").Value = 851.55817
$ws.Range("I122").Value = 660.63635
$ws.Range("J122").Value = 917.1875
$ws.Range("K122").Value = 5945.72715
$ws.Range("L122").Value = 8254.6875
$ws.Range("M122").Value = -3495.72715
$ws.Range("N122").Value = -13154.6875
# row 129
$ws.Range("H129").Value = 13889915
$ws.Range("I129").Value = 33333860
$ws.Range("J129").Value = 4167942.2
$ws.Range("K129").Value = 100001580
$ws.Range("L129").Value = 12503826.6
$ws.Range("M129").Value = -99996580
$ws.Range("N129").Value = -12513826.6
# row 131
$ws.Range("H131").Value = 29456908
$ws.Range("I131").Value = 200000500
$ws.Range("J131").Value = 52840.137
$ws.Range("K131").Value = 600001500
$ws.Range("L131").Value = 158520.411
$ws.Range("M131").Value = -599996460
$ws.Range("N131").Value = -168600.411

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 49
$ws.Range("H49").Value = 21000
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 21000
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 21000
$ws.Range("N49").Value = -21368
$ws.Range("M49").ClearContents()
# row 95
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
# row 97
$ws.Range("H97").Value = 760.35297
$ws.Range("I97").Value = 754.3333
$ws.Range("J97").Value = 805.5
$ws.Range("K97").Value = 754.3333
$ws.Range("L97").Value = 805.5
$ws.Range("M97").Value = -258.3333
$ws.Range("N97").Value = -1797.5
# row 132
$ws.Range("H132").Value = 4772.091
$ws.Range("I132").Value = 5198.8486
$ws.Range("J132").Value = 3491.818
$ws.Range("K132").Value = 15596.5458
$ws.Range("L132").Value = 10475.454
$ws.Range("M132").Value = -13066.5458
$ws.Range("N132").Value = -15535.454
# row 141
$ws.Range("H141").Value = 55050
$ws.Range("J141").Value = 55050
$ws.Range("L141").Value = 55050
$ws.Range("N141").Value = -65410

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 40
$ws.Range("H40").Value = 2459.2144
$ws.Range("I40").Value = 2402.2727
$ws.Range("K40").Value = 2402.2727
$ws.Range("M40").Value = -2266.2727
# row 61
$ws.Range("H61").Value = 1398.9231
$ws.Range("I61").Value = 1265.5
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 1265.5
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -1063.5
$ws.Range("N61").Value = -3404
# row 113
$ws.Range("H113").Value = 1398.9231
$ws.Range("I113").Value = 1265.5
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 1265.5
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 904.5
$ws.Range("N113").Value = -7340
# row 131
$ws.Range("H131").Value = 30000
$ws.Range("J131").Value = 30000
$ws.Range("L131").Value = 30000
$ws.Range("N131").Value = -40080
# row 136
$ws.Range("H136").Value = 1628
$ws.Range("I136").Value = 1542.0714
$ws.Range("K136").Value = 4626.2142
$ws.Range("M136").Value = -2076.2142

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 49
$ws.Range("H49").Value = 8062
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 8062
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 8062
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value = -8522
# row 100
$ws.Range("H100").Value = 367
$ws.Range("I100").Value = 376.57144
$ws.Range("K100").Value = 753.14288
$ws.Range("M100").Value = -212.14288
# row 132
$ws.Range("H132").Value = 3882.275
$ws.Range("I132").Value = 4041.9688
$ws.Range("J132").Value = 3243.5
$ws.Range("K132").Value = 12125.9064
$ws.Range("L132").Value = 9730.5
$ws.Range("M132").Value = -9595.9064
$ws.Range("N132").Value = -14790.5
# row 136
$ws.Range("H136").Value = 1607.6666
$ws.Range("I136").Value = 626.2174
$ws.Range("J136").Value = 2510.6
$ws.Range("K136").Value = 1878.6522
$ws.Range("L136").Value = 7531.799999999999
$ws.Range("M136").Value = 671.3478
$ws.Range("N136").Value = -12631.8
